$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("cantidad") values get " metros" appended to the existing quantity text.
$ws.Range("B2").Value = "121450 metros"
$ws.Range("B3").Value = "90202 metros"
$ws.Range("B4").Value = "67386 metros"
$ws.Range("B5").Value = "67383 metros"
$ws.Range("B6").Value = "32687 metros"
$ws.Range("B7").Value = "12050 metros"
$ws.Range("B8").Value = "7956 metros"

# Column C ("unidad") changes from "unidades" to "metros" for every data row.
$ws.Range("C2:C8").Value = "metros"
